$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: A1 empty, B1 = "Idade", C1 empty
$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Value = "Idade"
$ws.Range("C1").Style = "Normal"

# Rows 2-6: A = "Gbairl", B = "19" (stored as text), C empty
$rows = 2,3,4,5,6
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "Gbairl"

    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = "19"
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).Style = "Normal"
}
